$wb = $excel.ActiveWorkbook

$oldId = "24f3fbdb-43db-46bb-a4aa-3a64f07f679e"
$newId = "380c2aa3-0ffe-4e46-8fde-b1c31755a04d"

$oldZhXlf = "$oldId.02db8b92cf30802664081aa8dbe6dc337d4cbd24.zh-cn.xlf"
$newZhXlf = "$newId.93362ad7770b06ea06a2cbee2c343342ed55a0c8.zh-cn.xlf"
$oldDeXlf = "$oldId.02db8b92cf30802664081aa8dbe6dc337d4cbd24.de-de.xlf"
$newDeXlf = "$newId.93362ad7770b06ea06a2cbee2c343342ed55a0c8.de-de.xlf"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newId.md"
}
$wsOverview.Range("G2").Value = "2016-09-05 21:09:36"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-09-05 21:09:31"

foreach ($hl in @($wsZh.Hyperlinks)) {
    if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 9) {
        $hl.Delete()
    }
}
$wsZh.Range("I2").Value = "PLACEHOLDER_EMPTY"
$wsZh.Range("J2").Value = "PLACEHOLDER_EMPTY"
$wsZh.Cells.Replace("PLACEHOLDER_EMPTY", "")
$wsZh.Range("I2").Style = "Normal"

$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZh.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-09-05 21:09:36"

foreach ($hl in @($wsDe.Hyperlinks)) {
    if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 9) {
        $hl.Delete()
    }
}
$wsDe.Range("I2").Value = "PLACEHOLDER_EMPTY"
$wsDe.Range("J2").Value = "PLACEHOLDER_EMPTY"
$wsDe.Cells.Replace("PLACEHOLDER_EMPTY", "")
$wsDe.Range("I2").Style = "Normal"

$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDe.Columns.Item(10).ColumnWidth = 20.833333333333332
